$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 2416.1177
$ws.Range("I28").Value = 349.91666
$ws.Range("J28").Value = 7375
$ws.Range("K28").Value = 349.91666
$ws.Range("L28").Value = 7375
$ws.Range("M28").Value = 135.08334
$ws.Range("N28").Value = -8345
$ws.Range("H40").Value = 1723.1428
$ws.Range("I40").Value = 2142.6365
$ws.Range("K40").Value = 2142.6365
$ws.Range("M40").Value = -1967.6365
$ws.Range("H64").Value = 2387.5557
$ws.Range("I64").Value = 2177.6
$ws.Range("J64").Value = 2650
$ws.Range("K64").Value = 2177.6
$ws.Range("L64").Value = 2650
$ws.Range("M64").Value = -1929.6
$ws.Range("N64").Value = -3146
$ws.Range("H67").Value = 2387.5557
$ws.Range("I67").Value = 2177.6
$ws.Range("J67").Value = 2650
$ws.Range("K67").Value = 2177.6
$ws.Range("L67").Value = 2650
$ws.Range("M67").Value = -1319.6
$ws.Range("N67").Value = -4366
$ws.Range("H70").Value = 2293.2222
$ws.Range("I70").Value = 2100
$ws.Range("J70").Value = 2389.8333
$ws.Range("K70").Value = 6300
$ws.Range("L70").Value = 7169.499899999999
$ws.Range("M70").Value = -6030
$ws.Range("N70").Value = -7709.499899999999
$ws.Range("H73").Value = 2293.2222
$ws.Range("I73").Value = 2100
$ws.Range("J73").Value = 2389.8333
$ws.Range("K73").Value = 6300
$ws.Range("L73").Value = 7169.499899999999
$ws.Range("M73").Value = -5364
$ws.Range("N73").Value = -9041.499899999999
$ws.Range("H76").Value = 69620.2
$ws.Range("I76").Value = 74343.07000000001
$ws.Range("K76").Value = 74343.07000000001
$ws.Range("M76").Value = -74028.07000000001
$ws.Range("H79").Value = 69620.2
$ws.Range("I79").Value = 74343.07000000001
$ws.Range("K79").Value = 74343.07000000001
$ws.Range("M79").Value = -73251.07000000001
$ws.Range("H129").Value = 880773.5
$ws.Range("I129").Value = 537.8570999999999
$ws.Range("J129").Value = 1079536.4
$ws.Range("K129").Value = 1613.5713
$ws.Range("L129").Value = 3238609.2
$ws.Range("M129").Value = 3386.4287
$ws.Range("N129").Value = -3248609.2

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 29453.621
$ws.Range("I32").Value = 28509.592
$ws.Range("J32").Value = 32174.646
$ws.Range("K32").Value = 28509.592
$ws.Range("L32").Value = 32174.646
$ws.Range("M32").Value = -28222.592
$ws.Range("N32").Value = -32748.646
$ws.Range("H61").Value = 2123.7666
$ws.Range("I61").Value = 1743.55
$ws.Range("K61").Value = 1743.55
$ws.Range("M61").Value = -1531.55
$ws.Range("H74").Value = 802.7206
$ws.Range("I74").Value = 685.16364
$ws.Range("K74").Value = 685.16364
$ws.Range("M74").Value = 188.83636
$ws.Range("H77").Value = 802.7206
$ws.Range("I77").Value = 685.16364
$ws.Range("K77").Value = 3425.8182
$ws.Range("M77").Value = 942.1818000000003
$ws.Range("H110").Value = 3417.3076
$ws.Range("I110").Value = 3655.6365
$ws.Range("J110").Value = 2106.5
$ws.Range("K110").Value = 3655.6365
$ws.Range("L110").Value = 2106.5
$ws.Range("M110").Value = -1610.6365
$ws.Range("N110").Value = -6196.5
$ws.Range("H132").Value = 5450.627
$ws.Range("I132").Value = 6299.6045
$ws.Range("J132").Value = 3169
$ws.Range("K132").Value = 18898.8135
$ws.Range("L132").Value = 9507
$ws.Range("M132").Value = -16368.8135
$ws.Range("N132").Value = -14567
$ws.Range("H136").Value = 2123.7666
$ws.Range("I136").Value = 1743.55
$ws.Range("K136").Value = 5230.65
$ws.Range("M136").Value = -2680.65

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H40").Value = 30000
$ws.Range("J40").Value = 30000
$ws.Range("L40").Value = 30000
$ws.Range("N40").Value = -30530
$ws.Range("H96").Value = 3751.5
$ws.Range("I96").Value = 3751.5
$ws.Range("J96").Value = 0
$ws.Range("K96").Value = 3751.5
$ws.Range("L96").Value = 0
$ws.Range("M96").Value = -1005.5
$ws.Range("N96").Value = $null
$ws.Range("H99").Value = 876
$ws.Range("I99").Value = 872
$ws.Range("J99").Value = 900
$ws.Range("K99").Value = 872
$ws.Range("L99").Value = 900
$ws.Range("M99").Value = 626
$ws.Range("N99").Value = -3896
$ws.Range("H105").Value = 2012.28
$ws.Range("I105").Value = 1833.1428
$ws.Range("K105").Value = 1833.1428
$ws.Range("M105").Value = -86.14280000000008
$ws.Range("H107").Value = 1697.8889
$ws.Range("I107").Value = 1721
$ws.Range("J107").Value = 1669
$ws.Range("K107").Value = 1721
$ws.Range("L107").Value = 1669
$ws.Range("M107").Value = 199
$ws.Range("N107").Value = -5509
$ws.Range("H134").Value = 3654.5881
$ws.Range("I134").Value = 3779.1765
$ws.Range("K134").Value = 11337.5295
$ws.Range("M134").Value = -8802.529500000001

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 965.73334
$ws.Range("I16").Value = 902.625
$ws.Range("K16").Value = 902.625
$ws.Range("M16").Value = -615.625
$ws.Range("H62").Value = 47623016
$ws.Range("I62").Value = 4264.6665
$ws.Range("J62").Value = 83337080
$ws.Range("K62").Value = 4264.6665
$ws.Range("L62").Value = 83337080
$ws.Range("M62").Value = -3640.6665
$ws.Range("N62").Value = -83338328
$ws.Range("H65").Value = 47623016
$ws.Range("I65").Value = 4264.6665
$ws.Range("J65").Value = 83337080
$ws.Range("K65").Value = 21323.3325
$ws.Range("L65").Value = 416685400
$ws.Range("M65").Value = -18203.3325
$ws.Range("N65").Value = -416691640
$ws.Range("H113").Value = 965.73334
$ws.Range("I113").Value = 902.625
$ws.Range("K113").Value = 902.625
$ws.Range("M113").Value = 1267.375

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H86").Value = 2558.8235
$ws.Range("I86").Value = 2192.3076
$ws.Range("J86").Value = 3750
$ws.Range("K86").Value = 6576.9228
$ws.Range("L86").Value = 11250
$ws.Range("M86").Value = -5390.9228
$ws.Range("N86").Value = -13622
$ws.Range("H89").Value = 2558.8235
$ws.Range("I89").Value = 2192.3076
$ws.Range("J89").Value = 3750
$ws.Range("K89").Value = 19730.7684
$ws.Range("L89").Value = 33750
$ws.Range("M89").Value = -13802.7684
$ws.Range("N89").Value = -45606

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 17004336
$ws.Range("I70").Value = 22178142
$ws.Range("J70").Value = 4687
$ws.Range("K70").Value = 22178142
$ws.Range("L70").Value = 4687
$ws.Range("M70").Value = -22177872
$ws.Range("N70").Value = -5227
$ws.Range("H73").Value = 17004336
$ws.Range("I73").Value = 22178142
$ws.Range("J73").Value = 4687
$ws.Range("K73").Value = 22178142
$ws.Range("L73").Value = 4687
$ws.Range("M73").Value = -22177206
$ws.Range("N73").Value = -6559
$ws.Range("H126").Value = 4398.5713
$ws.Range("I126").Value = 3950
$ws.Range("J126").Value = 4996.6665
$ws.Range("K126").Value = 11850
$ws.Range("L126").Value = 14989.9995
$ws.Range("M126").Value = -9380
$ws.Range("N126").Value = -19929.9995
$ws.Range("H132").Value = 5755.3335
$ws.Range("I132").Value = 6354.304
$ws.Range("J132").Value = 3787.2856
$ws.Range("K132").Value = 19062.912
$ws.Range("L132").Value = 11361.8568
$ws.Range("M132").Value = -16532.912
$ws.Range("N132").Value = -16421.8568

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2739.92
$ws.Range("I122").Value = 2872.111
$ws.Range("J122").Value = 2400
$ws.Range("K122").Value = 8616.332999999999
$ws.Range("L122").Value = 7200
$ws.Range("M122").Value = -6166.332999999999
$ws.Range("N122").Value = -12100
